$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-26 Tuesday" "2024-11-27 Wednesday"

Replace-Text "923÷8=115, 3" "761÷4=190, 1"
Replace-Text "925÷2=462, 1" "976÷9=108, 4"
Replace-Text "797÷4=199, 1" "939÷5=187, 4"
Replace-Text "407÷5=81, 2" "197÷3=65, 2"
Replace-Text "422÷7=60, 2" "409÷4=102, 1"

Replace-Text "826÷7=118, 0" "740÷5=148, 0"
Replace-Text "813÷5=162, 3" "827÷8=103, 3"
Replace-Text "159÷8=19, 7" "183÷9=20, 3"
Replace-Text "700÷6=116, 4" "308÷5=61, 3"
Replace-Text "809÷9=89, 8" "875÷4=218, 3"

Replace-Text "825÷6=137, 3" "708÷6=118, 0"
Replace-Text "245÷9=27, 2" "394÷5=78, 4"
Replace-Text "686÷6=114, 2" "280÷4=70, 0"
Replace-Text "109÷2=54, 1" "306÷9=34, 0"
Replace-Text "642÷2=321, 0" "959÷3=319, 2"

Replace-Text "114÷4=28, 2" "833÷3=277, 2"
Replace-Text "499÷9=55, 4" "556÷8=69, 4"
Replace-Text "566÷3=188, 2" "315÷2=157, 1"
Replace-Text "899÷9=99, 8" "446÷5=89, 1"
Replace-Text "937÷6=156, 1" "214÷7=30, 4"

Replace-Text "670÷7=95, 5" "885÷6=147, 3"
Replace-Text "275÷2=137, 1" "926÷6=154, 2"
Replace-Text "808÷9=89, 7" "129÷2=64, 1"
Replace-Text "803÷6=133, 5" "249÷7=35, 4"
Replace-Text "726÷6=121, 0" "355÷9=39, 4"
